# 8.5.2.xlsx — add a new "2022" column (S) mirroring the existing
# "2021" column (R): same per-row formatting, new data values.
# Two of the rows (8 and 36) are section headers with no numeric value,
# they get a brand-new bold+italic style instead of a plain numeric one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122
$xlPasteFormats = -4122

function Set-YearCell([int]$row, $value) {
    $ws.Range("R$row").Copy()
    $ws.Range("S$row").PasteSpecial($xlPasteFormats)
    $ws.Range("S$row").Value = $value
}

# Header row: year label 2022 (same style as R4/2021)
Set-YearCell 4 2022

# Kyrgyz Republic totals / women / men
Set-YearCell 5 4.9000000000000004
Set-YearCell 6 6.1
Set-YearCell 7 4

# "by territory" section header (row 8) — blank cell, new bold+italic font
$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial($xlPasteFormats)
$ws.Range("S8").Font.Bold = $true
$ws.Range("S8").Font.Italic = $true

# Batken oblast
Set-YearCell 9 6.1
Set-YearCell 10 12.4
Set-YearCell 11 3.2

# Djalal-Abad oblast
Set-YearCell 12 10.8
Set-YearCell 13 14.6
Set-YearCell 14 8.5

# Ysyk-Kul oblast
Set-YearCell 15 5.5
Set-YearCell 16 7.1
Set-YearCell 17 4.4000000000000004

# Naryn oblast
Set-YearCell 18 5.8
Set-YearCell 19 11.6
Set-YearCell 20 3.1

# Osh oblast
Set-YearCell 21 1.5
Set-YearCell 22 2.2999999999999998
Set-YearCell 23 1

# Talas oblast
Set-YearCell 24 2.2999999999999998
Set-YearCell 25 3.3
Set-YearCell 26 1.6

# Chui oblast
Set-YearCell 27 4.5999999999999996
Set-YearCell 28 4.4000000000000004
Set-YearCell 29 4.7

# Bishkek city
Set-YearCell 30 4
Set-YearCell 31 3.2
Set-YearCell 32 4.7

# Osh city
Set-YearCell 33 2.6
Set-YearCell 34 3.3
Set-YearCell 35 2.2000000000000002

# "by age group" section header (row 36) — blank cell, new bold+italic font
$ws.Range("R36").Copy()
$ws.Range("S36").PasteSpecial($xlPasteFormats)
$ws.Range("S36").Font.Bold = $true
$ws.Range("S36").Font.Italic = $true

# Age groups 15-19 .. 60-69
Set-YearCell 37 13.2
Set-YearCell 38 7.5
Set-YearCell 39 4.0999999999999996
Set-YearCell 40 4.3
Set-YearCell 41 2.6
Set-YearCell 42 1

# 70 and over — footer row, value is the existing "…" shared string
$ws.Range("R43").Copy()
$ws.Range("S43").PasteSpecial($xlPasteFormats)
$ws.Range("S43").Value = "…"

# Clear clipboard marching ants / leftover copy mode
$excel.CutCopyMode = $false

# Match the recorded selection after the edit
$ws.Range("T12").Select() | Out-Null
